$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.290.47"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.71%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.298.69"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.51%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.97"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.80%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -8.44%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.300.50"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.45%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.61%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.19%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.86%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.863.59"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.64"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -6.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.298.10"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.62%  "

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.266.80"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.04"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -8.12%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.71"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.50"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.11%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.20%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.53"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -4.65%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.531"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -6.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.430.33"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.69%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -9.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.175"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.03"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.92%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.41"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "22.53"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.23"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -7.29%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.09"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.31%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.03"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -9.44%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.62"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.84%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.329.10"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0722"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -7.94%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "25.39"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -18.12%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.54%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.31%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.12"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -4.12%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.08"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -7.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.57"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -6.33%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.320.50"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -9.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.53"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -6.45%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -7.23%  "
